$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 730, pushing existing rows 730-771 down to 731-772.
$ws.Rows.Item(730).Insert()

# Populate the newly inserted row with the new data point.
# The date column stores plain text (e.g. "2026/01/31"), not a real date
# serial, so force text formatting before assigning the value - otherwise
# Excel auto-converts the slash-separated string into a date number.
$ws.Range("A730").NumberFormat = "@"
$ws.Range("A730").Value = "2026/01/31"
$ws.Range("A730").Style = "Normal"

$ws.Range("B730").Value = "土"
$ws.Range("C730").Value = 16
$ws.Range("D730").Value = 24
